$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (incl. date number format on column A) from the last
# existing data row (197) down onto the new rows (198:204) before writing
# values, so the new rows pick up the same cell styles (e.g. date style on A).
$ws.Range("A197:T197").Copy() | Out-Null
$ws.Range("A198:T204").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 198 ---
$ws.Range("A198").Value = 45815
$ws.Range("B198").Value = "Flowering"
$ws.Range("C198").Value = "Large"
$ws.Range("D198").Value = 64
$ws.Range("E198").Value = 77
$ws.Range("G198").Formula = "=1.37+0.62"
$ws.Range("H198").Value = 0.2
$ws.Range("I198").Value = "Yes"
$ws.Range("J198").Value = 2
$ws.Range("K198").Value = "Dark"
$ws.Range("L198").Value = 4
$ws.Range("M198").Value = 0.73
$ws.Range("N198").Value = 65
$ws.Range("O198").Value = 29.82
$ws.Range("P198").Value = 5
$ws.Range("Q198").Value = 0.83
$ws.Range("R198").Value = 9.3
$ws.Range("S198").Value = 55
$ws.Range("T198").Value = 23

# --- Row 199 ---
$ws.Range("A199").Value = 45815
$ws.Range("B199").Value = "Nonflowering"
$ws.Range("C199").Value = "Medium"
$ws.Range("D199").Value = 64
$ws.Range("E199").Value = 77
$ws.Range("H199").Value = 0.3
$ws.Range("I199").Value = "Yes"
$ws.Range("J199").Value = 3
$ws.Range("K199").Value = "Dark"
$ws.Range("L199").Value = 4
$ws.Range("M199").Value = 0.73
$ws.Range("N199").Value = 65
$ws.Range("O199").Value = 29.82
$ws.Range("P199").Value = 5
$ws.Range("Q199").Value = 0.83
$ws.Range("R199").Value = 9.3
$ws.Range("S199").Value = 55
$ws.Range("T199").Value = 23

# --- Row 200 ---
$ws.Range("A200").Value = 45815
$ws.Range("B200").Value = "Nonflowering"
$ws.Range("C200").Value = "Small"
$ws.Range("D200").Value = 64
$ws.Range("E200").Value = 77
$ws.Range("H200").Value = 1.5
$ws.Range("I200").Value = "Yes"
$ws.Range("J200").Value = 3
$ws.Range("K200").Value = "Dark"
$ws.Range("L200").Value = 4
$ws.Range("M200").Value = 0.73
$ws.Range("N200").Value = 65
$ws.Range("O200").Value = 29.82
$ws.Range("P200").Value = 5
$ws.Range("Q200").Value = 0.83
$ws.Range("R200").Value = 9.3
$ws.Range("S200").Value = 55
$ws.Range("T200").Value = 23

# --- Row 201 ---
$ws.Range("A201").Value = 45815
$ws.Range("B201").Value = "Nonflowering"
$ws.Range("C201").Value = "Medium"
$ws.Range("D201").Value = 64
$ws.Range("E201").Value = 77
$ws.Range("H201").Value = 1.25
$ws.Range("I201").Value = "Yes"
$ws.Range("J201").Value = 3
$ws.Range("K201").Value = "Neutral"
$ws.Range("L201").Value = 4
$ws.Range("M201").Value = 0.73
$ws.Range("N201").Value = 65
$ws.Range("O201").Value = 29.82
$ws.Range("P201").Value = 5
$ws.Range("Q201").Value = 0.83
$ws.Range("R201").Value = 9.3
$ws.Range("S201").Value = 55
$ws.Range("T201").Value = 23

# --- Row 202 ---
$ws.Range("A202").Value = 45815
$ws.Range("B202").Value = "Nonflowering"
$ws.Range("C202").Value = "Medium"
$ws.Range("D202").Value = 64
$ws.Range("E202").Value = 77
$ws.Range("H202").Value = 0.5
$ws.Range("I202").Value = "Yes"
$ws.Range("J202").Value = 3
$ws.Range("K202").Value = "Neutral"
$ws.Range("L202").Value = 4
$ws.Range("M202").Value = 0.73
$ws.Range("N202").Value = 65
$ws.Range("O202").Value = 29.82
$ws.Range("P202").Value = 5
$ws.Range("Q202").Value = 0.83
$ws.Range("R202").Value = 9.3
$ws.Range("S202").Value = 55
$ws.Range("T202").Value = 23

# --- Row 203 ---
$ws.Range("A203").Value = 45815
$ws.Range("B203").Value = "Nonflowering"
$ws.Range("C203").Value = "Large"
$ws.Range("D203").Value = 64
$ws.Range("E203").Value = 77
$ws.Range("H203").Value = 0.75
$ws.Range("I203").Value = "Yes"
$ws.Range("J203").Value = 4
$ws.Range("K203").Value = "Bright"
$ws.Range("L203").Value = 4
$ws.Range("M203").Value = 0.73
$ws.Range("N203").Value = 65
$ws.Range("O203").Value = 29.82
$ws.Range("P203").Value = 5
$ws.Range("Q203").Value = 0.83
$ws.Range("R203").Value = 9.3
$ws.Range("S203").Value = 55
$ws.Range("T203").Value = 23

# --- Row 204 ---
$ws.Range("A204").Value = 45815
$ws.Range("B204").Value = "Tree"
$ws.Range("C204").Value = "Medium"
$ws.Range("D204").Value = 64
$ws.Range("E204").Value = 77
$ws.Range("H204").Formula = "=4/3"
$ws.Range("I204").Value = "Yes"
$ws.Range("J204").Value = 1
$ws.Range("K204").Value = "Neutral"
$ws.Range("L204").Value = 4
$ws.Range("M204").Value = 0.73
$ws.Range("N204").Value = 65
$ws.Range("O204").Value = 29.82
$ws.Range("P204").Value = 5
$ws.Range("Q204").Value = 0.83
$ws.Range("R204").Value = 9.3
$ws.Range("S204").Value = 55
$ws.Range("T204").Value = 23

# Shared formulas: F195:F197 already share one ABS formula; extend the
# pattern down through the new rows (198:204) as a second shared block, and
# give column G (constant 1.37+0.62) the same treatment for rows 199:204
# (row 198 gets a standalone copy of the formula, matching the source file).
$ws.Range("F198:F204").Formula = "=ABS(D198-E198)"
$ws.Range("G198").Formula = "=1.37+0.62"
$ws.Range("G199:G204").Formula = "=1.37+0.62"

# Update the sheet view to match where the author left the selection after
# adding the rows.
$ws.Range("I200").Select() | Out-Null
